# Update VII_IND / VII_CNJ "Order" column values to their abbreviated
# forms, and adjust the active selections on each sheet, matching the
# "Update tags to align with lang-otw #1" commit.

$wb = $excel.ActiveWorkbook

$wsInd = $wb.Worksheets.Item("VII_IND")
$wsCnj = $wb.Worksheets.Item("VII_CNJ")

# VII_IND: column B holds "Independent" for every data row (2:257) -> "Ind"
$indRange = $wsInd.Range("B2:B257")
$indRange.Value = "Ind"

# VII_CNJ: column B holds "Conjunct" for every data row (2:129) -> "Cnj"
$cnjRange = $wsCnj.Range("B2:B129")
$cnjRange.Value = "Cnj"

# Update selections / views to match the saved state in the workbook.
$wsInd.Activate() | Out-Null
$wsInd.Range("B2:B257").Select() | Out-Null

$wsCnj.Activate() | Out-Null
$wsCnj.Range("H2").Select() | Out-Null

$wsInd.Activate() | Out-Null
